$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "=""29.220.14"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E2").Value = "  -0.62%  "

$c = $ws.Range("D3")
$c.Formula = "=""1.860.46"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E3").Value = "  -0.92%  "

$c = $ws.Range("D4")
$c.Formula = "=""0.9992"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.15%  "

$c = $ws.Range("D5")
$c.Formula = "=""242.33"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.63%  "

$c = $ws.Range("D6")
$c.Formula = "=""0.7027"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E6").Value = "  -1.99%  "

$c = $ws.Range("D7")
$c.Formula = "=""0.9994"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.14%  "

$c = $ws.Range("D8")
$c.Formula = "=""0.07828"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E8").Value = "  -1.67%  "

$c = $ws.Range("D9")
$c.Formula = "=""0.3122"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E9").Value = "  -0.83%  "

$c = $ws.Range("D10")
$c.Formula = "=""24.08"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E10").Value = "  -3.50%  "

$c = $ws.Range("D11")
$c.Formula = "=""0.07808"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E11").Value = "  -3.96%  "

$c = $ws.Range("D12")
$c.Formula = "=""1.864.14"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.91%  "

$c = $ws.Range("D13")
$c.Formula = "=""5.142"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E13").Value = "  -1.84%  "

$c = $ws.Range("D14")
$c.Formula = "=""92.59"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E14").Value = "  -2.87%  "

$c = $ws.Range("D15")
$c.Formula = "=""0.6936"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E15").Value = "  -1.95%  "

$c = $ws.Range("D16")
$c.Formula = "=""6.584"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.66%  "

$c = $ws.Range("D17")
$c.Formula = "=""0.000008516"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E17").Value = "  +0.97%  "

$c = $ws.Range("D18")
$c.Formula = "=""29.202.99"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.70%  "

$c = $ws.Range("D19")
$c.Formula = "=""250.11"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.24%  "

$c = $ws.Range("D20")
$c.Formula = "=""2.106.26"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E20").Value = "  -1.25%  "

$c = $ws.Range("D21")
$c.Formula = "=""12.97"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E21").Value = "  -3.05%  "

$c = $ws.Range("D22")
$c.Formula = "=""0.9992"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E22").Value = "  -0.16%  "

$c = $ws.Range("D23")
$c.Formula = "=""7.617"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E23").Value = "  -0.70%  "

$c = $ws.Range("D24")
$c.Formula = "=""0.9997"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.04%  "

$c = $ws.Range("D25")
$c.Formula = "=""0.1541"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E25").Value = "  -3.03%  "

$c = $ws.Range("D26")
$c.Formula = "=""160.85"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.71%  "

$c = $ws.Range("D27")
$c.Formula = "=""8.935"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E27").Value = "  -1.45%  "

$c = $ws.Range("D28")
$c.Formula = "=""18.63"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E28").Value = "  -1.52%  "

$c = $ws.Range("D29")
$c.Formula = "=""1.574"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E29").Value = "  +4.39%  "

$c = $ws.Range("D30")
$c.Formula = "=""4.288"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E30").Value = "  -2.98%  "

$c = $ws.Range("D31")
$c.Formula = "=""4.250"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E31").Value = "  -1.28%  "

$c = $ws.Range("D32")
$c.Formula = "=""1.207"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E32").Value = "  -0.97%  "

$c = $ws.Range("D33")
$c.Formula = "=""0.05251"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.43%  "

$c = $ws.Range("D34")
$c.Formula = "=""0.7611"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.63%  "

$c = $ws.Range("D35")
$c.Formula = "=""1.877"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E35").Value = "  -3.71%  "

$c = $ws.Range("D36")
$c.Formula = "=""1.179"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E36").Value = "  +0.27%  "

$c = $ws.Range("D37")
$c.Formula = "=""2.705"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E37").Value = "  +0.12%  "

$c = $ws.Range("D38")
$c.Formula = "=""0.01867"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E38").Value = "  -1.43%  "

$c = $ws.Range("D39")
$c.Formula = "=""1.237.23"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E39").Value = "  -2.85%  "

$c = $ws.Range("D40")
$c.Formula = "=""2.720"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E40").Value = "  -1.61%  "

$c = $ws.Range("D41")
$c.Formula = "=""0.9019"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.53%  "

$c = $ws.Range("D42")
$c.Formula = "=""110.25"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E42").Value = "  -1.69%  "

$c = $ws.Range("D43")
$c.Formula = "=""5.854"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E43").Value = "  -8.45%  "

$c = $ws.Range("D44")
$c.Formula = "=""0.9990"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.16%  "

$c = $ws.Range("D45")
$c.Formula = "=""68.51"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E45").Value = "  -7.78%  "

$c = $ws.Range("D46")
$c.Formula = "=""2.005.08"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("E47").Value = "  -4.33%  "

$c = $ws.Range("D48")
$c.Formula = "=""0.5183"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E48").Value = "  -0.33%  "

$c = $ws.Range("D49")
$c.Formula = "=""9.537"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.14%  "

$c = $ws.Range("D50")
$c.Formula = "=""1.768"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E50").Value = "  -2.17%  "

$c = $ws.Range("D51")
$c.Formula = "=""0.4263"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.91%  "

$excel.CutCopyMode = $false